$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val.EndsWith("16")) {
        $cell.Value = $val.Substring(0, $val.Length - 2)
    }
}
